$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: the "@<SUMM_TEXT>@" placeholder in the Borrower ("Заемщиком")
# clause was accidentally split across three separate runs
# ("@<SUMM_" / "TEXT" / ">@"). Re-merge it into a single run carrying
# the formatting of the first of the three runs (Times New Roman, bold).
# There is another, already-correct, single-run "@<SUMM_TEXT>@" earlier
# in the document (in the Guarantor clause), so anchor on the unique
# preceding sentence to make sure we edit the right one.
# ---------------------------------------------------------------------
$anchor = $d.Content
[void]$anchor.Find.Execute("Договор вступает в силу", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
$anchor.Collapse(0)
[void]$anchor.Find.Execute("@<SUMM_TEXT>@", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
# Setting .Text to the exact same string the range already holds is a
# no-op for the run layout, so bounce through a scratch value first to
# force the three runs to collapse into one before writing the real text.
$anchor.Text = "TEMP_PLACEHOLDER"
$anchor.Text = "@<SUMM_TEXT>@"

# ---------------------------------------------------------------------
# Change 2: drop the stray trailing spaces after the hard-coded BIK
# number "040349602" in the bank-details table.
# ---------------------------------------------------------------------
[void]$d.Content.Find.Execute("040349602            ", $false, $false, $false, $false, $false, `
                               $true, 1, $false, "040349602", 2)

# ---------------------------------------------------------------------
# Change 3: remove the "Корр. счет: @<CORR_ACCOUNT>@" line entirely
# (bold label run + value run + its trailing line break), keeping the
# line break that follows "@<BIK>@" so the paragraph still ends the
# same way.
# ---------------------------------------------------------------------
$bik = $d.Content
[void]$bik.Find.Execute("@<BIK>@", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0)
$afterBik = $bik.End

$corr = $d.Content
[void]$corr.Find.Execute("@<CORR_ACCOUNT>@", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)
$corrEnd = $corr.End + 1   # include the <w:br/> right after the CORR_ACCOUNT tag

$toRemove = $d.Range($afterBik + 1, $corrEnd)   # keep the <w:br/> right after @<BIK>@
$toRemove.Delete()
